$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 2168.5
$ws.Range("I99").Value = 2684.75
$ws.Range("J99").Value = 1480.1666
$ws.Range("K99").Value = 8054.25
$ws.Range("L99").Value = 4440.4998
$ws.Range("M99").Value = -6556.25
$ws.Range("N99").Value = -7436.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 48000
$ws.Range("J109").Value = 48000
$ws.Range("L109").Value = 48000
$ws.Range("N109").Value = -50774
$ws.Range("H132").Value = 18246
$ws.Range("I132").Value = 15494.5
$ws.Range("J132").Value = 22648.4
$ws.Range("K132").Value = 46483.5
$ws.Range("L132").Value = 67945.20000000001
$ws.Range("M132").Value = -43953.5
$ws.Range("N132").Value = -73005.20000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 210.9375
$ws.Range("I7").Value = 282.81818
$ws.Range("J7").Value = 52.8
$ws.Range("K7").Value = 282.81818
$ws.Range("L7").Value = 52.8
$ws.Range("M7").Value = -169.81818
$ws.Range("N7").Value = -278.8
$ws.Range("H22").Value = 958.25
$ws.Range("I22").Value = 105.25
$ws.Range("J22").Value = 1811.25
$ws.Range("K22").Value = 105.25
$ws.Range("L22").Value = 1811.25
$ws.Range("M22").Value = 244.75
$ws.Range("N22").Value = -2511.25
$ws.Range("H41").Value = 9180
$ws.Range("J41").Value = 11350
$ws.Range("L41").Value = 11350
$ws.Range("N41").Value = -12206
$ws.Range("H50").Value = 6769
$ws.Range("I50").Value = 2527.6667
$ws.Range("J50").Value = 9950
$ws.Range("K50").Value = 2527.6667
$ws.Range("L50").Value = 9950
$ws.Range("M50").Value = -1902.6667
$ws.Range("N50").Value = -11200
$ws.Range("H51").Value = 45000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 45000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 45000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -46472
$ws.Range("H58").Value = 1015.9804
$ws.Range("I58").Value = 798.24243
$ws.Range("J58").Value = 1415.1666
$ws.Range("K58").Value = 798.24243
$ws.Range("L58").Value = 1415.1666
$ws.Range("M58").Value = -595.24243
$ws.Range("N58").Value = -1821.1666
$ws.Range("H59").Value = 19750
$ws.Range("J59").Value = 19750
$ws.Range("L59").Value = 19750
$ws.Range("N59").Value = -22040
$ws.Range("H60").Value = 20100
$ws.Range("J60").Value = 20100
$ws.Range("L60").Value = 20100
$ws.Range("N60").Value = -21122
$ws.Range("H61").Value = 45000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 45000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 45000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -45696
$ws.Range("H74").Value = 20142.5
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 20142.5
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 3250.2666
$ws.Range("I86").Value = 2743.9375
$ws.Range("J86").Value = 3828.9285
$ws.Range("K86").Value = 2743.9375
$ws.Range("L86").Value = 3828.9285
$ws.Range("M86").Value = -1620.9375
$ws.Range("N86").Value = -6074.9285
$ws.Range("H89").Value = 3250.2666
$ws.Range("I89").Value = 2743.9375
$ws.Range("J89").Value = 3828.9285
$ws.Range("K89").Value = 13719.6875
$ws.Range("L89").Value = 19144.6425
$ws.Range("M89").Value = -8103.6875
$ws.Range("N89").Value = -30376.6425
$ws.Range("H136").Value = 1015.9804
$ws.Range("I136").Value = 798.24243
$ws.Range("J136").Value = 1415.1666
$ws.Range("K136").Value = 2394.72729
$ws.Range("L136").Value = 4245.4998
$ws.Range("M136").Value = 155.2727100000002
$ws.Range("N136").Value = -9345.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 6655.091
$ws.Range("I17").Value = 333.33334
$ws.Range("J17").Value = 9025.75
$ws.Range("K17").Value = 1000.00002
$ws.Range("L17").Value = 27077.25
$ws.Range("M17").Value = -831.0000200000001
$ws.Range("N17").Value = -27415.25
$ws.Range("H34").Value = 620.38464
$ws.Range("I34").Value = 192.4
$ws.Range("J34").Value = 887.875
$ws.Range("K34").Value = 577.2
$ws.Range("L34").Value = 2663.625
$ws.Range("M34").Value = -493.2
$ws.Range("N34").Value = -2831.625
$ws.Range("H60").Value = 160.42857
$ws.Range("I60").Value = 160.42857
$ws.Range("K60").Value = 481.28571
$ws.Range("M60").Value = -230.28571
$ws.Range("H131").Value = 774.6818
$ws.Range("I131").Value = 560.5909
$ws.Range("J131").Value = 988.7727
$ws.Range("K131").Value = 1681.7727
$ws.Range("L131").Value = 2966.3181
$ws.Range("M131").Value = 3358.2273
$ws.Range("N131").Value = -13046.3181
$ws.Range("H137").Value = 4440303
$ws.Range("I137").Value = 62273.832
$ws.Range("J137").Value = 20201208
$ws.Range("K137").Value = 186821.496
$ws.Range("L137").Value = 60603624
$ws.Range("M137").Value = -181721.496
$ws.Range("N137").Value = -60613824

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 35095.855
$ws.Range("J104").Value = 35095.855
$ws.Range("L104").Value = 35095.855
$ws.Range("N104").Value = -42083.855
$ws.Range("H123").Value = 22703.75
$ws.Range("J123").Value = 22703.75
$ws.Range("L123").Value = 22703.75
$ws.Range("N123").Value = -27603.75
$ws.Range("H126").Value = 1360.9
$ws.Range("I126").Value = 1244.5454
$ws.Range("J126").Value = 1503.1111
$ws.Range("K126").Value = 3733.6362
$ws.Range("L126").Value = 4509.3333
$ws.Range("M126").Value = -1263.6362
$ws.Range("N126").Value = -9449.3333
$ws.Range("H132").Value = 3717.0193
$ws.Range("I132").Value = 4926.879
$ws.Range("J132").Value = 1615.6842
$ws.Range("K132").Value = 14780.637
$ws.Range("L132").Value = 4847.0526
$ws.Range("M132").Value = -12250.637
$ws.Range("N132").Value = -9907.052599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 74605.57
$ws.Range("I122").Value = 145486.86
$ws.Range("J122").Value = 3724.2856
$ws.Range("K122").Value = 436460.58
$ws.Range("L122").Value = 11172.8568
$ws.Range("M122").Value = -434010.58
$ws.Range("N122").Value = -16072.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9092927
$ws.Range("I122").Value = 12501684
$ws.Range("J122").Value = 2906.6667
$ws.Range("K122").Value = 37505052
$ws.Range("L122").Value = 8720.000100000001
$ws.Range("M122").Value = -37502602
$ws.Range("N122").Value = -13620.0001
$ws.Range("H124").Value = 48000
$ws.Range("J124").Value = 48000
$ws.Range("L124").Value = 48000
$ws.Range("N124").Value = -57820
